# Split the sentence that announces the manuscript title so that the
# (new) title sits in its own run, matching:
#   "We are submitting our manuscript titled "<TITLE>" for consideration..."
#
# Old title:
#   Plant species identity matters when comparing the trophic impacts of
#   native and non-native plants: insights from a community-wide
#   bird-exclusion experiment
#
# New title:
#   Are native plants always superior foraging opportunities for
#   insectivores compared to invasives?

$d = $word.ActiveDocument

$quoteOpen  = [char]0x201C   # “
$quoteClose = [char]0x201D   # ”

$oldTitle = "Plant species identity matters when comparing the trophic " + `
            "impacts of native and non-native plants: insights from a " + `
            "community-wide bird-exclusion experiment"

$newTitle = "Are native plants always superior foraging opportunities " + `
            "for insectivores compared to invasives?"

# Locate the exact run of text that is the old manuscript title (the text
# strictly between the curly quotes), leaving the surrounding
# "We are submitting...titled " / " for consideration..." text untouched.
$titleRange = $d.Content
$found = $titleRange.Find.Execute($oldTitle, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the old manuscript title text to replace."
}
$titleStart = $titleRange.Start

# Replace the title text itself.
$titleRange.Text = $newTitle
$titleEnd = $titleStart + $newTitle.Length

# Re-address the (now new) title text by its known character offsets, then
# nudge a character-formatting property on it (set then immediately
# restore) purely to force Word to split the paragraph's run at the title's
# boundaries, so the title ends up isolated in its own run instead of being
# merged back into the surrounding text's run.
$titleRange = $d.Range($titleStart, $titleEnd)
$titleRange.Bold = 1
$titleRange.Bold = 0

Write-Output "Replaced manuscript title with new run split."
